# Update confusion-matrix cell counts on the "Matriz de Confusão" sheet
# to reflect the imported model's updated predictions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("S3").Value = 2
$ws.Range("X3").Value = 1
$ws.Range("E5").Value = 4
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 1
$ws.Range("X5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 4
$ws.Range("L7").Value = 1
$ws.Range("U7").Value = 0
$ws.Range("H8").Value = 4
$ws.Range("N8").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("X8").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("S9").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 7
$ws.Range("D12").Value = 0
$ws.Range("L12").Value = 3
$ws.Range("N12").Value = 0
$ws.Range("Q12").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("Y12").Value = 1
$ws.Range("H13").Value = 0
$ws.Range("T13").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("X14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("S15").Value = 2
$ws.Range("D16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("P16").Value = 5
$ws.Range("R16").Value = 2
$ws.Range("B17").Value = 0
$ws.Range("Q17").Value = 7
$ws.Range("H18").Value = 1
$ws.Range("N18").Value = 1
$ws.Range("R18").Value = 5
$ws.Range("S18").Value = 0
$ws.Range("B20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("S20").Value = 0
$ws.Range("T20").Value = 5
$ws.Range("E22").Value = 0
$ws.Range("N22").Value = 0
$ws.Range("V22").Value = 6
$ws.Range("X22").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("N23").Value = 0
$ws.Range("X23").Value = 0
$ws.Range("Z23").Value = 1
$ws.Range("H24").Value = 0
$ws.Range("N24").Value = 1
$ws.Range("S24").Value = 1
$ws.Range("X24").Value = 5
$ws.Range("B25").Value = 0
$ws.Range("Y25").Value = 6
$ws.Range("Y26").Value = 1
$ws.Range("Z26").Value = 5
$ws.Range("AA27").Value = 5
$ws.Range("B27").Value = 1
$ws.Range("K27").Value = 1
